$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:F4").NumberFormat = "@"

$ws.Range("A2").Value = "0.831"
$ws.Range("B2").Value = "0.794"
$ws.Range("C2").Value = "0.595"
$ws.Range("D2").Value = "0.887"
$ws.Range("E2").Value = "0.627"
$ws.Range("F2").Value = "0.978"

$ws.Range("A3").Value = "0.837"
$ws.Range("B3").Value = "0.824"
$ws.Range("C3").Value = "0.558"
$ws.Range("D3").Value = "0.908"
$ws.Range("E3").Value = "0.618"

$ws.Range("A4").Value = "0.841"
$ws.Range("B4").Value = "0.840"
$ws.Range("C4").Value = "0.588"
$ws.Range("D4").Value = "0.828"
$ws.Range("E4").Value = "0.629"
$ws.Range("F4").Value = "0.982"
